# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (period in arrears) rows for each trabajador are
# reordered so period 1710 precedes period 1711 for every worker, and
# the table is refreshed with the corresponding N° Doc / Nombre / Valor
# Mora / Salario Basico values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: WILMER PUELLO GRAU, period 1711 -> 1710 (Valor Mora / Salario Basico unchanged)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1050944927"
$ws.Range("D16").Value = "WILMER PUELLO GRAU"
$ws.Range("E16").Value = "1710"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 737717

# Row 17: now ANGI CARINA HERRERA SALAZAR, period 1710
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143341464"
$ws.Range("D17").Value = "ANGI CARINA HERRERA SALAZAR"
$ws.Range("E17").Value = "1710"
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1200000

# Row 18: now WILMER PUELLO GRAU, period 1711
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1050944927"
$ws.Range("D18").Value = "WILMER PUELLO GRAU"
$ws.Range("E18").Value = "1711"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 737717

# Row 19: ANGI CARINA HERRERA SALAZAR, period 1710 -> 1711
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143341464"
$ws.Range("D19").Value = "ANGI CARINA HERRERA SALAZAR"
$ws.Range("E19").Value = "1711"
$ws.Range("F19").Value = 48000
$ws.Range("G19").Value = 1200000
